$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest year (2009) by deleting row 2; this shifts all
# subsequent rows up by one (2010 -> row2, ..., 2020 -> row12).
$ws.Rows.Item(2).Delete()

# Append the newly reported year (2021) as the new last row (13).
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 2.4119
$ws.Range("C13").Value = 0.4037
$ws.Range("D13").Value = 0.7964
$ws.Range("E13").Value = 1.2201
$ws.Range("F13").Value = 1.4664
$ws.Range("G13").Value = 7.6655
$ws.Range("H13").Value = 20.7918
$ws.Range("I13").Value = 0.2777
$ws.Range("J13").Value = 3.7476
$ws.Range("K13").Value = 0.052
$ws.Range("L13").Value = 3.2482
$ws.Range("M13").Value = 1.9193
$ws.Range("N13").Value = 3.99
$ws.Range("O13").Value = 8.372400000000001
$ws.Range("P13").Value = 30.4883
$ws.Range("Q13").Value = 0.208
$ws.Range("R13").Value = 0.1645
$ws.Range("S13").Value = 23.7104
$ws.Range("T13").Value = $null
$ws.Range("U13").Value = 19.9496
$ws.Range("V13").Value = 1.9377

# Match the label-column styling (bold, centered, bordered) used by
# the other year cells in column A.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
